$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles/borders) from column K (row 3-5) into the new column L
$ws.Range("K3:K5").Copy()
$ws.Range("L3:L5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new column values: header year 2021 and its data value 269
$ws.Range("L4").Value = 2021
$ws.Range("L5").Value = 269

# Update the active selection shown in the saved view
$ws.Range("N3").Select()
